$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<india>"
$ws.Range("C2").Value = 56

$ws.Range("B3").Value = "<gape>"
$ws.Range("C3").Value = 50

$ws.Range("B4").Value = "<sene>"
$ws.Range("C4").Value = 54

$ws.Range("C5").Value = 53

$ws.Range("B6").Value = "<its>"
$ws.Range("C6").Value = 51

$ws.Range("B7").Value = "<whiskey>"
$ws.Range("C7").Value = 50

$ws.Range("B9").Value = "<it>"
$ws.Range("C9").Value = 12
